# Weekly fruit/vegetable price update.
# Insert two new rows (new weekly observations) right before the current
# row 369, pushing all existing data rows down by two (old 369..416
# becomes new 371..418), then populate the two freshly inserted rows
# with the new week's "Primera" / "Segunda" records.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 369:416 down to 371:418, leaving two blank rows at 369:370.
$ws.Rows("369:370").Insert()

# New row 369 - "Primera" quality record for the new week.
$ws.Range("A369").Value = 1
$ws.Range("B369").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C369").Value = "Arica y Parinacota"
$ws.Range("D369").Value = 44946
$ws.Range("E369").Value = 15
$ws.Range("F369").Value = 100112043
$ws.Range("G369").Value = "Pepino ensalada"
$ws.Range("H369").Value = "Sin especificar"
$ws.Range("I369").Value = "Primera"
$ws.Range("J369").Value = 160
$ws.Range("K369").Value = 6000
$ws.Range("L369").Value = 7000
$ws.Range("M369").Value = 6500
$ws.Range("N369").Value = "$/caja 70 unidades"
$ws.Range("O369").Value = "Región de Arica y Parinacota"
$ws.Range("P369").Value = 93
$ws.Range("Q369").Value = 70
$ws.Range("R369").Value = "Hortaliza"

# New row 370 - "Segunda" quality record for the new week.
$ws.Range("A370").Value = 1
$ws.Range("B370").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C370").Value = "Arica y Parinacota"
$ws.Range("D370").Value = 44946
$ws.Range("E370").Value = 15
$ws.Range("F370").Value = 100112043
$ws.Range("G370").Value = "Pepino ensalada"
$ws.Range("H370").Value = "Sin especificar"
$ws.Range("I370").Value = "Segunda"
$ws.Range("J370").Value = 160
$ws.Range("K370").Value = 5000
$ws.Range("L370").Value = 6000
$ws.Range("M370").Value = 5500
$ws.Range("N370").Value = "$/caja 100 unidades"
$ws.Range("O370").Value = "Región de Arica y Parinacota"
$ws.Range("P370").Value = 55
$ws.Range("Q370").Value = 100
$ws.Range("R370").Value = "Hortaliza"
